$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 15: Thu 17-Aug-2023 (serial 45155) ---
$ws.Range("A14").Copy()
$ws.Range("A15").PasteSpecial(-4122)  # xlPasteFormats - pick up the date style (s="4")
$ws.Range("A15").Value = 45155
$ws.Range("B15:F15").Value = "PRESENT"
$ws.Range("G15:J15").Value = "ABSENT"

# --- New row 16: Fri 18-Aug-2023 (serial 45156) ---
$ws.Range("A14").Copy()
$ws.Range("A16").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("A16").Value = 45156
$ws.Range("B16:F16").Value = "PRESENT"
$ws.Range("G16:J16").Value = "ABSENT"

# --- Comments for the newly-absent cells ---
$ws.Range("G15").AddComment("Dell:`nNot Informed")
$ws.Range("H15").AddComment("Dell:`nNot Informed")
$ws.Range("G16").AddComment("Dell:`nNot Informed")
$ws.Range("H16").AddComment("Dell:`nNot Informed")

# --- Extend the date-column validation to cover the new rows, keep validation order ---
$ws.Range("A2:A14").Validation.Delete()
$ws.Range("B2:K1048576").Validation.Delete()
$ws.Range("A2:A16").Validation.Add(0)
$ws.Range("A2:A16").Validation.IgnoreBlank = $false
$ws.Range("B2:K1048576").Validation.Add(3, 1, 1, '"PRESENT, ABSENT"') | Out-Null

# --- Move the active selection like the saved workbook did ---
$ws.Range("J16").Select() | Out-Null
